$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of C2, C3, and D3 as per the target diff
$ws.Range("C2").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
